# edit.ps1 - apply workshop.pptx update:
#  * slide1 date range update
#  * rework the "Notes" slide (slide 36) into "Part IV - Maintenance"
#  * add three new "Part IV.x" slides after it
#  * the original "Notes" slide content is preserved, relocated to the end of the deck

function Set-StructuredText($range, $paragraphs) {
    # Build the full run of text, with paragraphs joined by carriage returns.
    $full = ""
    $first = $true
    foreach ($para in $paragraphs) {
        if (-not $first) { $full += "`r" }
        $first = $false
        foreach ($run in $para.Runs) {
            $full += $run.Text
        }
    }
    $range.Text = $full

    # Paragraph-level formatting (bullet / numbering / indent level).
    for ($i = 1; $i -le $paragraphs.Count; $i++) {
        $para = $paragraphs[$i - 1]
        $pr = $range.Paragraphs($i)
        $pr.IndentLevel = 1
        if ($para.Bullet -eq "number") {
            $pr.ParagraphFormat.Bullet.Type = 2
        } else {
            $pr.ParagraphFormat.Bullet.Type = 0
        }
    }

    # Run-level character formatting (Courier font / italics) via absolute offsets.
    $pos = 1
    foreach ($para in $paragraphs) {
        foreach ($run in $para.Runs) {
            $len = $run.Text.Length
            if ($len -gt 0) {
                $sub = $range.Characters($pos, $len)
                if ($run.Courier) { $sub.Font.Name = "Courier" }
                if ($run.Italic) { $sub.Font.Italic = $true }
            }
            $pos = $pos + $len
        }
        $pos = $pos + 1
    }
}

function Set-Title($slide, $text) {
    $title = $slide.Shapes.Item(1)
    $tr = $title.TextFrame.TextRange
    $tr.Text = $text
    $pr = $tr.Paragraphs(1)
    $pr.IndentLevel = 1
    $pr.ParagraphFormat.Bullet.Type = 0
}

$p = $ppt.ActivePresentation

# 1. Update the workshop date range on the title slide.
$p.Slides.Item(1).Shapes.Item(3).TextFrame.TextRange.Text = "2025-08-26/2025-09-17"

# 2. Duplicate the current "Notes" slide (slide 36) so its content/hyperlink survives
#    unchanged; the duplicate will be pushed down to the end of the deck (slide 40)
#    once the new Part IV slides are inserted after the (soon to be retitled) slide 36.
$notesDup = $p.Slides.Item(36).Duplicate()

# 3. Rework slide 36 in place: "Notes" -> "Part IV - Maintenance".
$s36 = $p.Slides.Item(36)
Set-Title $s36 "Part IV - Maintenance"

$body36 = $s36.Shapes.Item(2).TextFrame.TextRange
$paras36 = @(
    @{ Bullet = "none"; Runs = @(
        @{ Text = "After initial deposit, metadata in Zotero records may be updated. Also, PDFs associated with a Zotero record may be exchanged with a different one (e.g., an incorrect pdf was associated with some literature record)."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "none"; Runs = @(
        @{ Text = "To update metadata of a Zenodo deposit associated with a Zotero record, you can re-run Step III.4 after including either "; Courier = $false; Italic = $false },
        @{ Text = "--update-metadata-only"; Courier = $true; Italic = $false },
        @{ Text = " or "; Courier = $false; Italic = $false },
        @{ Text = "--new-version"; Courier = $true; Italic = $false },
        @{ Text = "."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "none"; Runs = @(
        @{ Text = "With "; Courier = $false; Italic = $false },
        @{ Text = "--update-metadata-only"; Courier = $true; Italic = $false },
        @{ Text = " an metadata of an existing Zenodo record is updated with the Zotero record metadata. No new version is created and the pdf attachment is left untouched."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "none"; Runs = @(
        @{ Text = "With "; Courier = $false; Italic = $false },
        @{ Text = "--new-version"; Courier = $true; Italic = $false },
        @{ Text = " a new version of an existing Zenodo record is deposited with updated the Zotero record metadata and the associated pdf in Zotero."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "none"; Runs = @(
        @{ Text = "Note that Zenodo record metadata is editable, however Zenodo record files are "; Courier = $false; Italic = $false },
        @{ Text = "not"; Courier = $false; Italic = $true },
        @{ Text = "."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "none"; Runs = @(
        @{ Text = "So, when you need to update a pdf associated with a Zotero record, you need to create a "; Courier = $false; Italic = $false },
        @{ Text = "--new-version"; Courier = $true; Italic = $false },
        @{ Text = ". This is "; Courier = $false; Italic = $false },
        @{ Text = "not"; Courier = $false; Italic = $true },
        @{ Text = " done by default."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "none"; Runs = @(
        @{ Text = "Default behavior is to not update the associated Zenodo record and skip the deposit."; Courier = $false; Italic = $false }
    ) }
)
Set-StructuredText $body36 $paras36

# 4. Insert the three new "Part IV.x" slides right after slide 36 (positions 37-39),
#    using the same "Title and Content" layout as the rest of the deck.
$layout = 2

$s37 = $p.Slides.Add(37, $layout)
Set-Title $s37 "Part IV.1 - Edit Existing Record and Update Metadata"
$body37 = $s37.Shapes.Item(2).TextFrame.TextRange
$paras37 = @(
    @{ Bullet = "number"; Runs = @(
        @{ Text = "Update metadata for an already deposited Zotero record."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "number"; Runs = @(
        @{ Text = "Run the deposit workflow with the "; Courier = $false; Italic = $false },
        @{ Text = "--update-metadata-only"; Courier = $true; Italic = $false }
    ) },
    @{ Bullet = "number"; Runs = @(
        @{ Text = "Verify that the metadata of record in Zenodo was updated, but no new version was created."; Courier = $false; Italic = $false }
    ) }
)
Set-StructuredText $body37 $paras37

$s38 = $p.Slides.Add(38, $layout)
Set-Title $s38 "Part IV.2 - Create New Record with Updated PDF and Metadata"
$body38 = $s38.Shapes.Item(2).TextFrame.TextRange
$paras38 = @(
    @{ Bullet = "number"; Runs = @(
        @{ Text = "Replace a pdf attachment for an already deposited Zotero record."; Courier = $false; Italic = $false }
    ) },
    @{ Bullet = "number"; Runs = @(
        @{ Text = "Run the deposit workflow with the "; Courier = $false; Italic = $false },
        @{ Text = "--new-version"; Courier = $true; Italic = $false }
    ) },
    @{ Bullet = "number"; Runs = @(
        @{ Text = "Verify that a new version was created for the Zenodo record including the updated pdf"; Courier = $false; Italic = $false }
    ) }
)
Set-StructuredText $body38 $paras38

$s39 = $p.Slides.Add(39, $layout)
Set-Title $s39 "Part IV.3 - Retire Zenodo Deposit Associated with Deleted Zotero Record"
$body39 = $s39.Shapes.Item(2).TextFrame.TextRange
$paras39 = @(
    @{ Bullet = "none"; Runs = @(
        @{ Text = "(for now, manual workflow, can be automated if needed) 1. Locate the Zenodo Deposit Associated with a Zotero record that no longer exists 2. Click on Community :gear: icon “submit to community” in lower right panel 3. Submit to the “batlit-retired” community 4. Click on “manage communities” 5. If present, remove the deposit from the BatLit and BLR communities"; Courier = $false; Italic = $false }
    ) }
)
Set-StructuredText $body39 $paras39

Write-Host "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    $t = $sl.Shapes.Item(1).TextFrame.TextRange.Text
    Write-Host "$i`: $t"
}
